# Auto-generated Excel COM-interop edit script
# Applies 'Update latest output (run 188)' changes to optimisation_result workbook

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# ---- Sheet 'Schedule' updates (rows 2-5) ----
$wsSchedule.Range("B2").Value = 46073.22916666666
$wsSchedule.Range("C2").Value = 5.5
$wsSchedule.Range("D2").Value = 20.79
$wsSchedule.Range("E2").Value = 895.3878277499999
$wsSchedule.Range("F2").Value = 43.06819758297258
$wsSchedule.Range("A3").Value = 46073.3125
$wsSchedule.Range("C3").Value = 8.5
$wsSchedule.Range("D3").Value = 32.13
$wsSchedule.Range("E3").Value = 826.7582895000002
$wsSchedule.Range("F3").Value = 25.73166167133521
$wsSchedule.Range("A4").Value = 46073.91666666666
$wsSchedule.Range("B4").Value = 46074.08333333334
$wsSchedule.Range("E4").Value = 680.28356175
$wsSchedule.Range("F4").Value = 44.99229905753968
$wsSchedule.Range("A5").Value = 46074.3125
$wsSchedule.Range("B5").Value = 46074.72916666666
$wsSchedule.Range("C5").Value = 10
$wsSchedule.Range("D5").Value = 37.8
$wsSchedule.Range("E5").Value = 262.600767
$wsSchedule.Range("F5").Value = 6.947110238095239

# ---- Sheet 'Detailed' updates (rows 12, 17, 38-97) ----
$wsDetailed.Range("E12").Value = "ON"
$wsDetailed.Range("E17").Value = "ON"
$wsDetailed.Range("B38").Value = 91.68945
$wsDetailed.Range("B39").Value = 154.2
$wsDetailed.Range("B40").Value = 182.50982
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 299.75
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 252.98315
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 279.9547
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 169.78702
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("E44").Value = "OFF"
$wsDetailed.Range("B45").Value = 136.96511
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("E45").Value = "OFF"
$wsDetailed.Range("B46").Value = 108.89
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("B47").Value = 87.14085
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("B48").Value = 78.0
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("B49").Value = 84.79
$wsDetailed.Range("B50").Value = 93.67687
$wsDetailed.Range("B51").Value = 84.79
$wsDetailed.Range("B52").Value = 80.20633
$wsDetailed.Range("E52").Value = "ON"
$wsDetailed.Range("B53").Value = 80.23268
$wsDetailed.Range("E53").Value = "ON"
$wsDetailed.Range("B54").Value = 80.30043
$wsDetailed.Range("B55").Value = 79.95016
$wsDetailed.Range("B56").Value = 80.45996
$wsDetailed.Range("B57").Value = 80.50974
$wsDetailed.Range("B58").Value = 80.53968
$wsDetailed.Range("B59").Value = 80.80479
$wsDetailed.Range("B60").Value = 80.27793
$wsDetailed.Range("B61").Value = 92.87231
$wsDetailed.Range("B62").Value = 93.91743
$wsDetailed.Range("B63").Value = 95.45453
$wsDetailed.Range("B64").Value = 59.14226
$wsDetailed.Range("E64").Value = "OFF"
$wsDetailed.Range("B65").Value = 12.48745
$wsDetailed.Range("B66").Value = 1.6283
$wsDetailed.Range("B69").Value = 0.51
$wsDetailed.Range("B70").Value = 0.01049
$wsDetailed.Range("B71").Value = 0.51
$wsDetailed.Range("B72").Value = 0.36381
$wsDetailed.Range("B75").Value = 0.36344
$wsDetailed.Range("B76").Value = 0.01096
$wsDetailed.Range("B78").Value = 12.15582
$wsDetailed.Range("B79").Value = 35.86
$wsDetailed.Range("B80").Value = 35.88
$wsDetailed.Range("B81").Value = 35.88
$wsDetailed.Range("B84").Value = 55.15385
$wsDetailed.Range("B85").Value = 57.31011
$wsDetailed.Range("E85").Value = "OFF"
$wsDetailed.Range("B86").Value = 84.79
$wsDetailed.Range("B88").Value = 134.10369
$wsDetailed.Range("B89").Value = 108.01
$wsDetailed.Range("B90").Value = 144.68466
$wsDetailed.Range("B92").Value = 98.77869
$wsDetailed.Range("B93").Value = 85.95
$wsDetailed.Range("B94").Value = 74.99234
$wsDetailed.Range("B95").Value = 78.0
$wsDetailed.Range("B96").Value = 78.0
$wsDetailed.Range("B97").Value = 78.0
